# filter() method for featureGroupsSet and negate support for features method
#
# Adds a new "done" (column G) marker to several rows of the fGroups
# implementation-status sheet, moves/changes a few "X"/"X?" markers, and
# annotates two rows with a "maybe wait for autoID branch" note in the new
# column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value constants used throughout the sheet ---
$X    = "X"
$NOTE = "maybe wait for autoID branch"

# --- New "done" (column G) markers -----------------------------------
$ws.Range("G7").Value  = $X
$ws.Range("G8").Value  = $X
$ws.Range("G12").Value = $X
$ws.Range("G13").Value = $X
$ws.Range("G14").Value = $X
$ws.Range("G24").Value = $X
$ws.Range("G30").Value = $X
$ws.Range("G31").Value = $X
$ws.Range("G32").Value = $X
$ws.Range("G34").Value = $X
$ws.Range("G35").Value = $X
$ws.Range("G36").Value = $X
$ws.Range("G37").Value = $X
$ws.Range("G38").Value = $X
$ws.Range("G44").Value = $X
$ws.Range("G50").Value = $X
$ws.Range("G51").Value = $X
$ws.Range("G53").Value = $X

# --- Row 8: "X?" -> "X" (now implemented) ------------------------------
$ws.Range("B8").Value = $X

# --- Row 13: the "X" for filter() moves from column B to column C ------
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = $X

# --- Row 24: the "X" moves from column C to column G --------------------
$ws.Range("C24").ClearContents()

# --- Rows whose "X?" became a confirmed "X" ------------------------------
$ws.Range("B30").Value = $X
$ws.Range("B31").Value = $X
$ws.Range("B34").Value = $X
$ws.Range("B44").Value = $X
$ws.Range("B50").Value = $X
$ws.Range("B51").Value = $X

# --- New notes in column H ------------------------------------------------
$ws.Range("H23").Value = $NOTE
$ws.Range("H48").Value = $NOTE

# --- Update current selection to reflect where the author left off -------
$ws.Range("G15").Select()
